$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.415.26'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.101.43'
$ws.Range("E3").Value = '  -0.17%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '334.43'
$ws.Range("E5").Value = '  +1.35%  '

$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5225'
$ws.Range("E7").Value = '  -0.75%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4553'
$ws.Range("E8").Value = '  +3.59%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.52'
$ws.Range("E9").Value = '  +14.59%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08896'
$ws.Range("E10").Value = '  +0.35%  '

$ws.Range("E11").Value = '  +1.24%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.12'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.087.07'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.795'
$ws.Range("E14").Value = '  +0.82%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.027'
$ws.Range("E15").Value = '  +3.31%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '97.06'
$ws.Range("E16").Value = '  +0.66%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001144'
$ws.Range("E17").Value = '  +1.23%  '

$ws.Range("E18").Value = '  -0.08%  '

$ws.Range("E19").Value = '  -0.31%  '

$ws.Range("E20").Value = '  +0.77%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  -0.01%  '

$ws.Range("E22").Value = '  -0.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.477.18'
$ws.Range("E23").Value = '  -0.26%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.354'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.333.90'
$ws.Range("E26").Value = '  -0.74%  '

$ws.Range("E27").Value = '  -1.19%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.37'
$ws.Range("E28").Value = '  +0.46%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.518'
$ws.Range("E29").Value = '  -3.18%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.01'
$ws.Range("E30").Value = '  +0.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.208'
$ws.Range("E31").Value = '  -0.19%  '

$ws.Range("E32").Value = '  -0.59%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.654'
$ws.Range("E33").Value = '  -1.41%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.402'
$ws.Range("E34").Value = '  +2.66%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.930'
$ws.Range("E35").Value = '  -0.09%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.39'
$ws.Range("E36").Value = '  +1.93%  '

$ws.Range("E37").Value = '  +6.45%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02573'
$ws.Range("E38").Value = '  -0.40%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06846'
$ws.Range("E39").Value = '  +1.97%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2314'
$ws.Range("E40").Value = '  +1.34%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.68'
$ws.Range("E41").Value = '  -0.45%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6875'
$ws.Range("E42").Value = '  +0.21%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.246'
$ws.Range("E43").Value = '  -1.81%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.321'
$ws.Range("E44").Value = '  +4.76%  '

$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6400'
$ws.Range("E45").Value = '  +0.10%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '14.00'
$ws.Range("E46").Value = '  -0.25%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.652'
$ws.Range("E47").Value = '  +0.54%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.247'
$ws.Range("E48").Value = '  -0.52%  '

$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '83.18'
$ws.Range("E49").Value = '  +1.09%  '

$ws.Range("B50").Value = 'WEMIXTOKEN'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.202'
$ws.Range("E50").Value = '  -1.01%  '

$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.00000000337'
$ws.Range("E51").Value = '  +13.83%  '
